$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 223.33333
$ws.Range("I2").Value = 223.33333
$ws.Range("K2").Value = 223.33333
$ws.Range("M2").Value = -110.33333

$ws.Range("H55").Value = 562.2963
$ws.Range("I55").Value = 545.5
$ws.Range("J55").Value = 610.2857
$ws.Range("K55").Value = 545.5
$ws.Range("L55").Value = 610.2857
$ws.Range("M55").Value = -331.5
$ws.Range("N55").Value = -1038.2857

$ws.Range("H62").Value = 5069.8
$ws.Range("I62").Value = 4837.25
$ws.Range("K62").Value = 4837.25
$ws.Range("M62").Value = -4213.25

$ws.Range("H65").Value = 5069.8
$ws.Range("I65").Value = 4837.25
$ws.Range("K65").Value = 24186.25
$ws.Range("M65").Value = -21066.25

$ws.Range("H98").Value = 0
$ws.Range("I98").Value = 0
$ws.Range("K98").Value = 0
$ws.Range("M98").ClearContents()

$ws.Range("H100").Value = 3839.75
$ws.Range("I100").Value = 1788
$ws.Range("K100").Value = 1788
$ws.Range("M100").Value = -1247

$ws.Range("H116").Value = 3174.5789
$ws.Range("I116").Value = 3040.3076
$ws.Range("J116").Value = 3465.5
$ws.Range("K116").Value = 3040.3076
$ws.Range("L116").Value = 3465.5
$ws.Range("M116").Value = 401.6923999999999
$ws.Range("N116").Value = -10349.5

$ws.Range("H122").Value = 0
$ws.Range("I122").Value = 0
$ws.Range("K122").Value = 0
$ws.Range("M122").ClearContents()

$ws.Range("H132").Value = 12425.909
$ws.Range("I132").Value = 1104.919
$ws.Range("K132").Value = 3314.757000000001
$ws.Range("M132").Value = -784.7570000000005

$ws.Range("H137").Value = 2449.9033
$ws.Range("I137").Value = 2270.2593
$ws.Range("K137").Value = 6810.777900000001
$ws.Range("M137").Value = -4260.777900000001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 3835
$ws.Range("I45").Value = 2840.2144
$ws.Range("K45").Value = 2840.2144
$ws.Range("M45").Value = -2463.2144

$ws.Range("H61").Value = 1798.2759
$ws.Range("I61").Value = 917.55
$ws.Range("K61").Value = 917.55
$ws.Range("M61").Value = -705.55

$ws.Range("H74").Value = 1310.742
$ws.Range("I74").Value = 1284.4667
$ws.Range("K74").Value = 1284.4667
$ws.Range("M74").Value = -410.4666999999999

$ws.Range("H77").Value = 1310.742
$ws.Range("I77").Value = 1284.4667
$ws.Range("K77").Value = 6422.3335
$ws.Range("M77").Value = -2054.3335

$ws.Range("H122").Value = 3562.9756
$ws.Range("I122").Value = 1676.7188
$ws.Range("K122").Value = 5030.1564
$ws.Range("M122").Value = -2580.1564

$ws.Range("H132").Value = 2381.7344
$ws.Range("I132").Value = 1929.1608
$ws.Range("K132").Value = 5787.482400000001
$ws.Range("M132").Value = -3257.482400000001

$ws.Range("H136").Value = 1798.2759
$ws.Range("I136").Value = 917.55
$ws.Range("K136").Value = 2752.65
$ws.Range("M136").Value = -202.6499999999996

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H12").Value = 201.33333
$ws.Range("I12").Value = 152
$ws.Range("J12").Value = 300
$ws.Range("K12").Value = 152
$ws.Range("L12").Value = 300
$ws.Range("M12").Value = 16
$ws.Range("N12").Value = -636

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H5").Value = 368.5
$ws.Range("J5").Value = 400
$ws.Range("L5").Value = 400
$ws.Range("N5").Value = -624

$ws.Range("H86").Value = 36718.77
$ws.Range("I86").Value = 42609.75
$ws.Range("J86").Value = 27293.2
$ws.Range("K86").Value = 42609.75
$ws.Range("L86").Value = 27293.2
$ws.Range("M86").Value = -41486.75
$ws.Range("N86").Value = -29539.2

$ws.Range("H89").Value = 36718.77
$ws.Range("I89").Value = 42609.75
$ws.Range("J89").Value = 27293.2
$ws.Range("K89").Value = 213048.75
$ws.Range("L89").Value = 136466
$ws.Range("M89").Value = -207432.75
$ws.Range("N89").Value = -147698

$ws.Range("H132").Value = 3277.4849
$ws.Range("I132").Value = 2337.5925
$ws.Range("K132").Value = 7012.7775
$ws.Range("M132").Value = -4482.7775

$ws.Range("H134").Value = 2893.9268
$ws.Range("I134").Value = 1509.3572
$ws.Range("J134").Value = 5876.077
$ws.Range("K134").Value = 4528.071599999999
$ws.Range("L134").Value = 17628.231
$ws.Range("M134").Value = -1993.071599999999
$ws.Range("N134").Value = -22698.231

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H7").Value = 1411.8462
$ws.Range("I7").Value = 1586.5454
$ws.Range("K7").Value = 4759.6362
$ws.Range("M7").Value = -4647.6362

$ws.Range("H39").Value = 3765
$ws.Range("J39").Value = 4207
$ws.Range("L39").Value = 12621
$ws.Range("N39").Value = -13209

$ws.Range("H55").Value = 23812844
$ws.Range("I55").Value = 228
$ws.Range("J55").Value = 30307194
$ws.Range("K55").Value = 684
$ws.Range("L55").Value = 90921582
$ws.Range("M55").Value = -507
$ws.Range("N55").Value = -90921936

$ws.Range("H107").Value = 298.22726
$ws.Range("I107").Value = 389.66666
$ws.Range("J107").Value = 283.78946
$ws.Range("K107").Value = 1168.99998
$ws.Range("L107").Value = 851.3683800000001
$ws.Range("M107").Value = 751.0000199999999
$ws.Range("N107").Value = -4691.36838

$ws.Range("H139").Value = 37039830
$ws.Range("I139").Value = 40002340
$ws.Range("K139").Value = 120007020
$ws.Range("M139").Value = -120001880

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H17").Value = 2092.3076
$ws.Range("J17").Value = 2092.3076
$ws.Range("L17").Value = 2092.3076
$ws.Range("N17").Value = -2428.3076

$ws.Range("H70").Value = 79914.336
$ws.Range("I70").Value = 127978.336
$ws.Range("K70").Value = 127978.336
$ws.Range("M70").Value = -127708.336

$ws.Range("H73").Value = 79914.336
$ws.Range("I73").Value = 127978.336
$ws.Range("K73").Value = 127978.336
$ws.Range("M73").Value = -127042.336

$ws.Range("H97").Value = 497.35
$ws.Range("I97").Value = 448
$ws.Range("J97").Value = 546.7
$ws.Range("K97").Value = 448
$ws.Range("L97").Value = 546.7
$ws.Range("M97").Value = 48
$ws.Range("N97").Value = -1538.7

$ws.Range("H122").Value = 5477.6
$ws.Range("I122").Value = 5477.6
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 16432.8
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -13982.8
$ws.Range("N122").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 8903.444
$ws.Range("I40").Value = 11019.583
$ws.Range("J40").Value = 7210.533
$ws.Range("K40").Value = 11019.583
$ws.Range("L40").Value = 7210.533
$ws.Range("M40").Value = -10883.583
$ws.Range("N40").Value = -7482.533

$ws.Range("H100").Value = 140420
$ws.Range("I100").Value = 278315.25
$ws.Range("J100").Value = 2524.75
$ws.Range("K100").Value = 278315.25
$ws.Range("L100").Value = 2524.75
$ws.Range("M100").Value = -277774.25
$ws.Range("N100").Value = -3606.75

$ws.Range("H132").Value = 3465.1365
$ws.Range("I132").Value = 2790.6562
$ws.Range("J132").Value = 5263.75
$ws.Range("K132").Value = 8371.9686
$ws.Range("L132").Value = 15791.25
$ws.Range("M132").Value = -5841.9686
$ws.Range("N132").Value = -20851.25

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 2076.7163
$ws.Range("I132").Value = 1782.2167
$ws.Range("K132").Value = 5346.6501
$ws.Range("M132").Value = -2816.6501
